# The deck ships two embedded theme parts:
#   ppt/theme/theme1.xml -> clrScheme "Office"   (the stock Office Theme colours)
#   ppt/theme/theme2.xml -> clrScheme "Integral" (the colour scheme actually
#                                                   wired to the slide master /
#                                                   presentation, i.e. the one
#                                                   that is "live" for every
#                                                   slide in the deck)
#
# The commit swaps the palette that is applied to the presentation: the
# slide master (and therefore every slide) stops using the "Integral"
# palette and starts using the stock "Office Theme" palette instead.
#
# The PowerPoint object model doesn't give scripts a way to rename a theme
# or colour scheme (ThemeColorScheme.Name is read-only, and there is no
# supported Save/Load-by-name for it), but it does let a script rewrite the
# twelve RGB slots of the theme colour scheme that is in effect for the
# presentation. Doing that reproduces the actual colour change described by
# the diff: every dk/lt/accent/hyperlink slot that used to hold the
# "Integral" values now holds the matching "Office Theme" values.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Office Theme colours (RGB packed as r + g*256 + b*65536, i.e. the classic
# VBA RGB() encoding used by ColorFormat.RGB / ThemeColor.RGB):
$cs.Item(1).RGB  = 0          # dk1      000000
$cs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Item(3).RGB  = 6968388    # dk2      44546A
$cs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Item(6).RGB  = 3243501    # accent2  ED7D31
$cs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Item(8).RGB  = 49407      # accent4  FFC000
$cs.Item(9).RGB  = 12874308   # accent5  4472C4
$cs.Item(10).RGB = 4697456    # accent6  70AD47
$cs.Item(11).RGB = 12673797   # hlink    0563C1
$cs.Item(12).RGB = 7491477    # folHlink 954F72
